# Corecion a Diebold Mariano y revision de Cap1
# Update the Diebold-Mariano summary table values:
#  - EnCQR-LSTM row: Comparaciones_Significativas "5/10" -> "4/10"; Proporcion_Sig 128 -> 102.4
#  - MCPS row: Comparaciones_Significativas "2/10" -> "1/10"; Proporcion_Sig 51.2 -> 25.6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "4/10"
$ws.Range("C2").Value = 102.4

$ws.Range("B4").Value = "1/10"
$ws.Range("C4").Value = 25.6
